$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2 through 72 gets updated to the new
# date serial value 45184 (2023-09-15).
for ($row = 2; $row -le 72; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
